# Updated cryptos list (Price / Volume(1h) refresh, plus two row swaps where
# the scraped ranking order changed) - mirrors the GitHub Actions data pull.
#
# NOTE: several "Price" values are numeric-looking text (e.g. "14.50",
# "2.20") that must stay stored as TEXT, matching the source data. Plain
# `Range.Value = '14.50'` would get auto-coerced to the Number 14.5 by
# Excel (silently dropping the trailing zero), so those literals are
# entered with a leading apostrophe ('' inside a single-quoted PowerShell
# string = one literal quote char), which is Excel's standard "force text"
# entry prefix and is stripped from the stored cell text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.837.04'
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('D3').Value = '2.210.72'
$ws.Range('E3').Value = '  -1.85%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '''257.51'
$ws.Range('E5').Value = '  +5.14%  '
$ws.Range('D6').Value = '''0.612'
$ws.Range('E6').Value = '  -0.99%  '
$ws.Range('D7').Value = '''76.89'
$ws.Range('E7').Value = '  +1.41%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -2.49%  '
$ws.Range('D10').Value = '''42.12'
$ws.Range('E10').Value = '  +2.39%  '
$ws.Range('D11').Value = '''0.0904'
$ws.Range('E11').Value = '  -3.40%  '
$ws.Range('D12').Value = '''6.96'
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('D14').Value = '2.543.49'
$ws.Range('E14').Value = '  -1.78%  '
$ws.Range('D15').Value = '''14.50'
$ws.Range('E15').Value = '  -1.35%  '
$ws.Range('D16').Value = '2.218.56'
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').Value = '''0.784'
$ws.Range('E17').Value = '  -2.53%  '
$ws.Range('D18').Value = '42.826.62'
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('D19').Value = '''0.0000103'
$ws.Range('E19').Value = '  -2.70%  '
$ws.Range('D20').Value = '''71.31'
$ws.Range('E20').Value = '  -0.23%  '
$ws.Range('D21').Value = '''5.97'
$ws.Range('E21').Value = '  -0.69%  '
$ws.Range('B22').Value = 'ImmutableX'
$ws.Range('C22').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D22').Value = '''2.20'
$ws.Range('E22').Value = '  +0.80%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Value = '''230.31'
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').Value = '''9.37'
$ws.Range('E24').Value = '  -6.76%  '
$ws.Range('D26').Value = '''42.42'
$ws.Range('E26').Value = '  +8.73%  '
$ws.Range('D27').Value = '''10.76'
$ws.Range('E27').Value = '  -1.68%  '
$ws.Range('D28').Value = '''3.35'
$ws.Range('E28').Value = '  -4.27%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = '''2.21'
$ws.Range('E29').Value = '  -2.04%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '''2.22'
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('D31').Value = '''173.83'
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('D32').Value = '''20.32'
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('E33').Value = '  +9.31%  '
$ws.Range('D34').Value = '''5.22'
$ws.Range('E34').Value = '  -2.05%  '
$ws.Range('E35').Value = '  -0.64%  '
$ws.Range('D36').Value = '''0.0358'
$ws.Range('E36').Value = '  +6.37%  '
$ws.Range('E37').Value = '  -3.24%  '
$ws.Range('E38').Value = '  -1.11%  '
$ws.Range('D39').Value = '''12.88'
$ws.Range('E39').Value = '  -2.01%  '
$ws.Range('E40').Value = '  +18.37%  '
$ws.Range('E41').Value = '  -1.33%  '
$ws.Range('E42').Value = '  -2.33%  '
$ws.Range('D43').Value = '''5.30'
$ws.Range('E43').Value = '  -4.15%  '
$ws.Range('D44').Value = '''60.08'
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').Value = '''102.66'
$ws.Range('E45').Value = '  -3.29%  '
$ws.Range('E46').Value = '  -4.78%  '
$ws.Range('D47').Value = '''0.0978'
$ws.Range('E47').Value = '  -1.89%  '
$ws.Range('D48').Value = '''0.463'
$ws.Range('E48').Value = '  -3.98%  '
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('E50').Value = '  -1.47%  '
$ws.Range('D51').Value = '2.430.47'
$ws.Range('E51').Value = '  -1.30%  '
